$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text (E1/F1) and hashtag row text (E2/F2) to reflect the
# "New Displacements" -> "Internal Displacements" terminology change.
$ws.Range("E1").Value = "Conflict Internal Displacements"
$ws.Range("F1").Value = "Disaster Internal Displacements"
$ws.Range("E2").Value = "#affected+idps+ind+internaldisp+conflict"
$ws.Range("F2").Value = "#affected+idps+ind+internaldisp+disaster"

# Move the active selection to F3 (matches the refreshed QuickCharts view).
$ws.Range("F3").Select()
